# Update cryptos list with latest prices / 1h volume changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "63.565.29"
$ws.Range("E2").Value = "  -6.40%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.256.47"
$ws.Range("E3").Value = "  -9.40%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.36%  "

# Row 5: Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "177.65"
$ws.Range("E5").Value = "  -12.14%  "

# Row 6: BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "512.43"
$ws.Range("E6").Value = "  -10.21%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.583"
$ws.Range("E7").Value = "  -5.35%  "

# Row 8: USDC (was LidoStakedEther)
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.08%  "

# Row 9: LidoStakedEther (was USDC)
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.257.22"
$ws.Range("E9").Value = "  -9.30%  "

# Row 10: Cardano
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.605"
$ws.Range("E10").Value = "  -11.09%  "

# Row 11: Avalanche
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.04"
$ws.Range("E11").Value = "  -5.81%  "

# Row 12: Dogecoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.128"
$ws.Range("E12").Value = "  -13.03%  "

# Row 13: ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  -11.98%  "

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.88"
$ws.Range("E14").Value = "  -13.82%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.803.01"
$ws.Range("E15").Value = "  -8.45%  "

# Row 16: TRON
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.120"
$ws.Range("E16").Value = "  -4.81%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "3.271.86"
$ws.Range("E17").Value = "  -8.93%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "63.459.38"
$ws.Range("E18").Value = "  -6.24%  "

# Row 19: Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.88"
$ws.Range("E19").Value = "  -11.78%  "

# Row 20: Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.60"
$ws.Range("E20").Value = "  -13.63%  "

# Row 21: Polygon
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.928"
$ws.Range("E21").Value = "  -12.71%  "

# Row 22: BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "364.68"
$ws.Range("E22").Value = "  -9.82%  "

# Row 23: Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "78.47"
$ws.Range("E23").Value = "  -7.44%  "

# Row 24: PancakeSwap
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.59"
$ws.Range("E24").Value = "  -14.38%  "

# Row 25: RenderToken
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.64"
$ws.Range("E25").Value = "  -16.37%  "

# Row 26: LEO
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.96"
$ws.Range("E26").Value = "  -2.83%  "

# Row 27: Toncoin
$ws.Range("E27").Value = "  -5.11%  "

# Row 28: ImmutableX
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.58"
$ws.Range("E28").Value = "  -10.92%  "

# Row 29: InternetComputer(DFINITY)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.93"
$ws.Range("E29").Value = "  -12.25%  "

# Row 30: Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.11"
$ws.Range("E30").Value = "  -12.06%  "

# Row 31: Bittensor
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "644.09"
$ws.Range("E31").Value = "  -4.75%  "

# Row 32: EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.89"
$ws.Range("E32").Value = "  -11.58%  "

# Row 33: NEARProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.52"
$ws.Range("E33").Value = "  -15.16%  "

# Row 34: Cosmos
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.87"
$ws.Range("E34").Value = "  -10.40%  "

# Row 35: OKB
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.89"
$ws.Range("E35").Value = "  -8.54%  "

# Row 36: Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.101"
$ws.Range("E36").Value = "  -11.03%  "

# Row 37: Dai
$ws.Range("E37").Value = "  +0.00%  "

# Row 38: InjectiveProtocol
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.30"
$ws.Range("E38").Value = "  -14.76%  "

# Row 39: TheGraph
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.364"
$ws.Range("E39").Value = "  -11.15%  "

# Row 40: FirstDigitalUSD
$ws.Range("E40").Value = "  +0.24%  "

# Row 41: Kaspa
$ws.Range("E41").Value = "  -9.79%  "

# Row 42: Maker
$ws.Range("D42").Value = "2.768.52"
$ws.Range("E42").Value = "  -13.50%  "

# Row 43: ThetaToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  -18.37%  "

# Row 44: WEMIXToken (was PEPE)
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  -8.22%  "

# Row 45: PEPE (was WEMIXToken)
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0612"
$ws.Range("E45").Value = "  -19.91%  "

# Row 46: VeChain
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0377"
$ws.Range("E46").Value = "  -8.07%  "

# Row 47: Fetch.AI
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.24"
$ws.Range("E47").Value = "  -17.07%  "

# Row 48: Stellar
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.121"
$ws.Range("E48").Value = "  -7.43%  "

# Row 49: Monero
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.66"
$ws.Range("E49").Value = "  -4.18%  "

# Row 50: Stacks
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.58"
$ws.Range("E50").Value = "  -4.64%  "

# Row 51: ApeXProtocol
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.78"
$ws.Range("E51").Value = "  -9.22%  "

